# chore: update Sheets via scheduled runner
# Refreshes market-price derived profit figures (currentAveragePrice* / Leve
# price & profit columns) across the ALC, ARM, BSM, CUL, GSM and WVR
# worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 53673.145
$ws.Range("J17").Value = 53673.145
$ws.Range("L17").Value = 161019.435
$ws.Range("N17").Value = -161355.435

# Row 52
$ws.Range("H52").Value = 251500
$ws.Range("J52").Value = 251500
$ws.Range("L52").Value = 754500
$ws.Range("N52").Value = -754820

# Row 64
$ws.Range("H64").Value = 3272.3076
$ws.Range("I64").Value = 3362.5
$ws.Range("J64").Value = 3128
$ws.Range("K64").Value = 3362.5
$ws.Range("L64").Value = 3128
$ws.Range("M64").Value = -3114.5
$ws.Range("N64").Value = -3624

# Row 67
$ws.Range("H67").Value = 3272.3076
$ws.Range("I67").Value = 3362.5
$ws.Range("J67").Value = 3128
$ws.Range("K67").Value = 3362.5
$ws.Range("L67").Value = 3128
$ws.Range("M67").Value = -2504.5
$ws.Range("N67").Value = -4844

# Row 87
$ws.Range("H87").Value = 26661.584
$ws.Range("J87").Value = 26661.584
$ws.Range("L87").Value = 26661.584
$ws.Range("N87").Value = -29157.584

# Row 90
$ws.Range("H90").Value = 26661.584
$ws.Range("J90").Value = 26661.584
$ws.Range("L90").Value = 79984.75199999999
$ws.Range("N90").Value = -92464.75199999999

# Row 111
$ws.Range("H111").Value = 766.1
$ws.Range("I111").Value = 604.8333
$ws.Range("J111").Value = 1008
$ws.Range("K111").Value = 1814.4999
$ws.Range("L111").Value = 3024
$ws.Range("M111").Value = 1252.5001
$ws.Range("N111").Value = -9158

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 14600
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 19400
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 19400
$ws.Range("M23").Value = -4741
$ws.Range("N23").Value = -19918

# Row 36
$ws.Range("H36").Value = 4202.8335
$ws.Range("I36").Value = 3643.4
$ws.Range("J36").Value = 7000
$ws.Range("K36").Value = 3643.4
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = -3297.4
$ws.Range("N36").Value = -7692

# Row 45
$ws.Range("H45").Value = 1502.25
$ws.Range("I45").Value = 1082.6207
$ws.Range("K45").Value = 1082.6207
$ws.Range("M45").Value = -705.6206999999999

# Row 63
$ws.Range("H63").Value = 2002.4
$ws.Range("I63").Value = 1780.4445
$ws.Range("K63").Value = 1780.4445
$ws.Range("M63").Value = -1094.4445

# Row 66
$ws.Range("H66").Value = 2002.4
$ws.Range("I66").Value = 1780.4445
$ws.Range("K66").Value = 8902.2225
$ws.Range("M66").Value = -5470.2225

# Row 74
$ws.Range("H74").Value = 1287.08
$ws.Range("I74").Value = 1170.3334
$ws.Range("J74").Value = 1900
$ws.Range("K74").Value = 1170.3334
$ws.Range("L74").Value = 1900
$ws.Range("M74").Value = -296.3334
$ws.Range("N74").Value = -3648

# Row 77
$ws.Range("H77").Value = 1287.08
$ws.Range("I77").Value = 1170.3334
$ws.Range("J77").Value = 1900
$ws.Range("K77").Value = 5851.666999999999
$ws.Range("L77").Value = 9500
$ws.Range("M77").Value = -1483.666999999999
$ws.Range("N77").Value = -18236

# Row 80
$ws.Range("H80").Value = 25941.8
$ws.Range("J80").Value = 25941.8
$ws.Range("L80").Value = 25941.8
$ws.Range("N80").Value = -27937.8

# Row 83
$ws.Range("H83").Value = 25941.8
$ws.Range("J83").Value = 25941.8
$ws.Range("L83").Value = 77825.39999999999
$ws.Range("N83").Value = -87809.39999999999

# Row 132
$ws.Range("H132").Value = 2956.3809
$ws.Range("I132").Value = 2940.3547
$ws.Range("J132").Value = 3001.5454
$ws.Range("K132").Value = 8821.0641
$ws.Range("L132").Value = 9004.636200000001
$ws.Range("M132").Value = -6291.0641
$ws.Range("N132").Value = -14064.6362

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 16571.666
$ws.Range("I35").Value = 4567
$ws.Range("J35").Value = 22574
$ws.Range("K35").Value = 4567
$ws.Range("L35").Value = 22574
$ws.Range("M35").Value = -4257
$ws.Range("N35").Value = -23194

# Row 68
$ws.Range("H68").Value = 25268
$ws.Range("I68").Value = 25268
$ws.Range("K68").Value = 25268
$ws.Range("M68").Value = -24457

# Row 71
$ws.Range("H71").Value = 25268
$ws.Range("I71").Value = 25268
$ws.Range("K71").Value = 75804
$ws.Range("M71").Value = -71748

# Row 99
$ws.Range("H99").Value = 3886
$ws.Range("I99").Value = 3724.5
$ws.Range("J99").Value = 3957.7778
$ws.Range("K99").Value = 3724.5
$ws.Range("L99").Value = 3957.7778
$ws.Range("M99").Value = -2226.5
$ws.Range("N99").Value = -6953.7778

# Row 102
$ws.Range("H102").Value = 18820
$ws.Range("I102").Value = 11500
$ws.Range("K102").Value = 11500
$ws.Range("M102").Value = -8255

# Row 107
$ws.Range("H107").Value = 4715.8
$ws.Range("I107").Value = 1399.5
$ws.Range("J107").Value = 6926.6665
$ws.Range("K107").Value = 1399.5
$ws.Range("L107").Value = 6926.6665
$ws.Range("M107").Value = 520.5
$ws.Range("N107").Value = -10766.6665

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 823.8889
$ws.Range("I113").Value = 399.375
$ws.Range("J113").Value = 1002.6316
$ws.Range("K113").Value = 1198.125
$ws.Range("L113").Value = 3007.8948
$ws.Range("M113").Value = 971.875
$ws.Range("N113").Value = -7347.8948

# Row 121
$ws.Range("H121").Value = 770.8182
$ws.Range("I121").Value = 226.33333
$ws.Range("J121").Value = 975
$ws.Range("K121").Value = 678.99999
$ws.Range("L121").Value = 2925
$ws.Range("M121").Value = 631.00001
$ws.Range("N121").Value = -5545

# Row 129
$ws.Range("H129").Value = 38772.5
$ws.Range("I129").Value = 5103
$ws.Range("K129").Value = 15309
$ws.Range("M129").Value = -10309

# Row 131
$ws.Range("H131").Value = 1333.6617
$ws.Range("J131").Value = 1107.4833
$ws.Range("L131").Value = 3322.449900000001
$ws.Range("N131").Value = -13402.4499

# Row 137
$ws.Range("H137").Value = 2492.8572
$ws.Range("I137").Value = 2394.2727
$ws.Range("J137").Value = 2854.3333
$ws.Range("K137").Value = 7182.8181
$ws.Range("L137").Value = 8562.999899999999
$ws.Range("M137").Value = -2082.8181
$ws.Range("N137").Value = -18762.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2379.3845
$ws.Range("I113").Value = 1804.5714
$ws.Range("J113").Value = 3050
$ws.Range("K113").Value = 1804.5714
$ws.Range("L113").Value = 3050
$ws.Range("M113").Value = 365.4286
$ws.Range("N113").Value = -7390

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1840.4286
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 1500
$ws.Range("M107").Value = 420
